# Fruta / hortaliza, semanal
# Insert two new weekly rows of "Durazno" (peach) price data for
# "Vega Modelo de Temuco" at the top of the existing data block (row 356),
# pushing all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 356-357; everything from the old row 356 onward
# shifts down by two rows (old 430/431 become new 432/433).
$ws.Rows("356:357").Insert()

# --- New row 356: Elegant Lady / Especial ---
$ws.Cells.Item(356, 1).Value = 10
$ws.Cells.Item(356, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(356, 3).Value = "La Araucanía"
$ws.Cells.Item(356, 4).Value = 45015
$ws.Cells.Item(356, 5).Value = 9
$ws.Cells.Item(356, 6).Value = "Fruta"
$ws.Cells.Item(356, 7).Value = 100103
$ws.Cells.Item(356, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(356, 9).Value = 100103004
$ws.Cells.Item(356, 10).Value = "Durazno"
$ws.Cells.Item(356, 11).Value = "Elegant Lady"
$ws.Cells.Item(356, 12).Value = "Especial"
$ws.Cells.Item(356, 13).Value = 125
$ws.Cells.Item(356, 14).Value = 26000
$ws.Cells.Item(356, 15).Value = 26000
$ws.Cells.Item(356, 16).Value = 26000
$ws.Cells.Item(356, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(356, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(356, 19).Value = 1444
$ws.Cells.Item(356, 20).Value = 18

# --- New row 357: Elegant Lady / Primera ---
$ws.Cells.Item(357, 1).Value = 10
$ws.Cells.Item(357, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(357, 3).Value = "La Araucanía"
$ws.Cells.Item(357, 4).Value = 45015
$ws.Cells.Item(357, 5).Value = 9
$ws.Cells.Item(357, 6).Value = "Fruta"
$ws.Cells.Item(357, 7).Value = 100103
$ws.Cells.Item(357, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(357, 9).Value = 100103004
$ws.Cells.Item(357, 10).Value = "Durazno"
$ws.Cells.Item(357, 11).Value = "Elegant Lady"
$ws.Cells.Item(357, 12).Value = "Primera"
$ws.Cells.Item(357, 13).Value = 125
$ws.Cells.Item(357, 14).Value = 22000
$ws.Cells.Item(357, 15).Value = 22000
$ws.Cells.Item(357, 16).Value = 22000
$ws.Cells.Item(357, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(357, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(357, 19).Value = 1222
$ws.Cells.Item(357, 20).Value = 18

# Make sure the date cells carry the same date number format as the rest
# of column D (style index 2 in the original workbook).
$ws.Cells.Item(356, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(357, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
